$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save", reusing the same formatting (style) as the
# other header cells (e.g. G1: bold, bordered, centered).
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data cells H2 / H3 with numeric value 1 (no special style).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
